$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value2 = 1.186522666666667
$ws.Cells.Item(2, 8).Value2 = 3.559568
$ws.Cells.Item(2, 9).Value2 = 0.06400371352898657
$ws.Cells.Item(2, 10).Value2 = 0.06400371352898658
$ws.Cells.Item(2, 11).Value2 = 3
$ws.Cells.Item(2, 12).Value2 = 1
$ws.Cells.Item(2, 13).Value2 = 0.3252056666666667
$ws.Cells.Item(2, 14).Value2 = 0.975617
$ws.Cells.Item(2, 15).Value2 = 0.0158278498560244
$ws.Cells.Item(2, 16).Value2 = 0.0158278498560244
$ws.Cells.Item(2, 17).Value2 = 0.3858638948284444
$ws.Cells.Item(2, 18).Value2 = 3.472775053456
$ws.Cells.Item(2, 19).Value2 = 0.001013041167964797
$ws.Cells.Item(2, 20).Value2 = 0.001013041167964797

# Row 3
$ws.Cells.Item(3, 7).Value2 = 1.186522666666667
$ws.Cells.Item(3, 8).Value2 = 3.559568
$ws.Cells.Item(3, 9).Value2 = 0.06400371352898657
$ws.Cells.Item(3, 10).Value2 = 0.06400371352898658
$ws.Cells.Item(3, 15).Value2 = 0.8133441666880411
$ws.Cells.Item(3, 16).Value2 = 0.8133441666880411
$ws.Cells.Item(3, 17).Value2 = 19.82835008223111
$ws.Cells.Item(3, 18).Value2 = 178.45515074008
$ws.Cells.Item(3, 19).Value2 = 0.05205704704517368
$ws.Cells.Item(3, 20).Value2 = 0.05205704704517369

# Row 4
$ws.Cells.Item(4, 7).Value2 = 1.186522666666667
$ws.Cells.Item(4, 8).Value2 = 3.559568
$ws.Cells.Item(4, 9).Value2 = 0.06400371352898657
$ws.Cells.Item(4, 10).Value2 = 0.06400371352898658
$ws.Cells.Item(4, 13).Value2 = 3.509903666666667
$ws.Cells.Item(4, 14).Value2 = 10.529711
$ws.Cells.Item(4, 15).Value2 = 0.1708279834559346
$ws.Cells.Item(4, 16).Value2 = 0.1708279834559346
$ws.Cells.Item(4, 17).Value2 = 4.164580258316445
$ws.Cells.Item(4, 18).Value2 = 37.48122232484801
$ws.Cells.Item(4, 19).Value2 = 0.01093362531584809
$ws.Cells.Item(4, 20).Value2 = 0.0109336253158481

# Row 5
$ws.Cells.Item(5, 9).Value2 = 0.599039184070822
$ws.Cells.Item(5, 10).Value2 = 0.599039184070822
$ws.Cells.Item(5, 11).Value2 = 3
$ws.Cells.Item(5, 12).Value2 = 1
$ws.Cells.Item(5, 13).Value2 = 0.3252056666666667
$ws.Cells.Item(5, 14).Value2 = 0.975617
$ws.Cells.Item(5, 15).Value2 = 0.0158278498560244
$ws.Cells.Item(5, 16).Value2 = 0.0158278498560244
$ws.Cells.Item(5, 17).Value2 = 3.611471584625111
$ws.Cells.Item(5, 18).Value2 = 32.503244261626
$ws.Cells.Item(5, 19).Value2 = 0.009481502263348334
$ws.Cells.Item(5, 20).Value2 = 0.009481502263348334

# Row 6
$ws.Cells.Item(6, 9).Value2 = 0.599039184070822
$ws.Cells.Item(6, 10).Value2 = 0.599039184070822
$ws.Cells.Item(6, 15).Value2 = 0.8133441666880411
$ws.Cells.Item(6, 16).Value2 = 0.8133441666880411
$ws.Cells.Item(6, 19).Value2 = 0.4872250259815668
$ws.Cells.Item(6, 20).Value2 = 0.4872250259815668

# Row 7
$ws.Cells.Item(7, 9).Value2 = 0.599039184070822
$ws.Cells.Item(7, 10).Value2 = 0.599039184070822
$ws.Cells.Item(7, 13).Value2 = 3.509903666666667
$ws.Cells.Item(7, 14).Value2 = 10.529711
$ws.Cells.Item(7, 15).Value2 = 0.1708279834559346
$ws.Cells.Item(7, 16).Value2 = 0.1708279834559346
$ws.Cells.Item(7, 17).Value2 = 38.97815645977312
$ws.Cells.Item(7, 18).Value2 = 350.8034081379581
$ws.Cells.Item(7, 19).Value2 = 0.1023326558259069
$ws.Cells.Item(7, 20).Value2 = 0.1023326558259069

# Row 8
$ws.Cells.Item(8, 7).Value2 = 6.246625666666667
$ws.Cells.Item(8, 8).Value2 = 18.739877
$ws.Cells.Item(8, 9).Value2 = 0.3369571024001913
$ws.Cells.Item(8, 10).Value2 = 0.3369571024001914
$ws.Cells.Item(8, 11).Value2 = 3
$ws.Cells.Item(8, 12).Value2 = 1
$ws.Cells.Item(8, 13).Value2 = 0.3252056666666667
$ws.Cells.Item(8, 14).Value2 = 0.975617
$ws.Cells.Item(8, 15).Value2 = 0.0158278498560244
$ws.Cells.Item(8, 16).Value2 = 0.0158278498560244
$ws.Cells.Item(8, 17).Value2 = 2.031438064345445
$ws.Cells.Item(8, 18).Value2 = 18.282942579109
$ws.Cells.Item(8, 19).Value2 = 0.005333306424711267
$ws.Cells.Item(8, 20).Value2 = 0.005333306424711268

# Row 9
$ws.Cells.Item(9, 7).Value2 = 6.246625666666667
$ws.Cells.Item(9, 8).Value2 = 18.739877
$ws.Cells.Item(9, 9).Value2 = 0.3369571024001913
$ws.Cells.Item(9, 10).Value2 = 0.3369571024001914
$ws.Cells.Item(9, 15).Value2 = 0.8133441666880411
$ws.Cells.Item(9, 16).Value2 = 0.8133441666880411
$ws.Cells.Item(9, 17).Value2 = 104.3893083806661
$ws.Cells.Item(9, 18).Value2 = 939.5037754259951
$ws.Cells.Item(9, 19).Value2 = 0.2740620936613006
$ws.Cells.Item(9, 20).Value2 = 0.2740620936613006

# Row 10
$ws.Cells.Item(10, 7).Value2 = 6.246625666666667
$ws.Cells.Item(10, 8).Value2 = 18.739877
$ws.Cells.Item(10, 9).Value2 = 0.3369571024001913
$ws.Cells.Item(10, 10).Value2 = 0.3369571024001914
$ws.Cells.Item(10, 13).Value2 = 3.509903666666667
$ws.Cells.Item(10, 14).Value2 = 10.529711
$ws.Cells.Item(10, 15).Value2 = 0.1708279834559346
$ws.Cells.Item(10, 16).Value2 = 0.1708279834559346
$ws.Cells.Item(10, 17).Value2 = 21.92505433172744
$ws.Cells.Item(10, 18).Value2 = 197.325488985547
$ws.Cells.Item(10, 19).Value2 = 0.05756170231417954
$ws.Cells.Item(10, 20).Value2 = 0.05756170231417955
